$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.472.89'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '3.094.14'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '622.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  -5.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.368'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '3.092.04'
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.739'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.198'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000252'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('D16').Value = '90.069.67'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '3.659.84'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.99%  '
$ws.Range('D19').Value = '3.090.34'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000220'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.59'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '437.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '89.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.160'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.200'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.12%  '
$ws.Range('E35').Value = '  +8.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '25.82'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.86'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.89%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.71%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '504.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.38%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0893'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.81%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.407'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +54.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.695'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '152.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('E51').Value = '  -0.10%  '
